$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.445.97'
$ws.Range("E2").Value = '  +2.19%  '
$ws.Range("D3").Value = '2.940.29'
$ws.Range("E3").Value = '  +1.87%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''588.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").Value = '''146.18'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.53%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '2.939.81'
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("D9").Value = '''0.507'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("D10").Value = '''7.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("E11").Value = '  +9.19%  '
$ws.Range("E12").Value = '  +1.20%  '
$ws.Range("D13").Value = '''0.0000233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.75%  '
$ws.Range("D14").Value = '''32.16'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '''0.125'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '3.429.47'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").Value = '62.424.74'
$ws.Range("E17").Value = '  +2.27%  '
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("D19").Value = '2.945.17'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").Value = '''433.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").Value = '''13.42'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").Value = '''0.661'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").Value = '''6.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("E24").Value = '  +5.21%  '
$ws.Range("D25").Value = '''80.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = '''11.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.22%  '
$ws.Range("D27").Value = '''2.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '''7.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.32%  '
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").Value = '''2.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.06%  '
$ws.Range("D32").Value = '0.0₃0999'
$ws.Range("E32").Value = '  +17.45%  '
$ws.Range("E33").Value = '  +3.63%  '
$ws.Range("D34").Value = '''26.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.90%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '''0.988'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.19%  '
$ws.Range("D37").Value = '''5.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("D38").Value = '''2.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.92%  '
$ws.Range("D39").Value = '''49.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.55%  '
$ws.Range("D40").Value = '''2.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.42%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("E43").Value = '  +3.43%  '
$ws.Range("D44").Value = '''38.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = '''135.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("D46").Value = '2.687.21'
$ws.Range("E46").Value = '  +1.05%  '
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("D48").Value = '''353.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.28%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("D51").Value = '''22.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.25%  '
